# bug fix in Eduati data files
# Mirrors the authored change: Sheet1 had 43 stray leftover rows (45-87,
# single column A) left over from a previous edit; trim them away so the
# sheet's used range matches the real data block (A1:N44). Also flip which
# sheet/cell is active: Sheet1 becomes the active tab (instead of Sheet3),
# with the cursor sitting mid-sheet, and Sheet3 is no longer the selected
# tab.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws3 = $wb.Worksheets.Item("Sheet3")

# --- Sheet1: drop the leftover rows 45:87 (only column A had values, a
#     remnant of stray data well past the real 44-row table) ---
$ws1.Range("A45:A87").EntireRow.Delete()

# --- Make Sheet1 the active sheet/tab (was Sheet3) ---
$ws1.Activate()
$ws1.Range("I49").Select()

# Sheet3 keeps its original selection (A2:N44 / active cell A2); only its
# tabSelected flag changes as a side effect of Sheet1 becoming active, so
# nothing further needs to be done here.

Write-Output "Trimmed Sheet1 to A1:N44 and switched the active tab to Sheet1."
